$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart Sprint5")

$ws.Range("I6").Value = 0.5
$ws.Range("I7").Value = 0.5
$ws.Range("H8").Value = 0.5
$ws.Range("I8").Value = 0.4
$ws.Range("F11").Value = 0.6
$ws.Range("G11").Value = 0.6
$ws.Range("H11").Value = 0.5
$ws.Range("I11").Value = 0.2
